$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 415
$ws.Range("I2").Value = 1144
$ws.Range("J2").Value = 4648
$ws.Range("K2").Value = 34
$ws.Range("L2").Value = 1246
$ws.Range("M2").Value = 74
$ws.Range("N2").Value = 834
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 22
$ws.Range("Q2").Value = 6
$ws.Range("R2").Value = 49
$ws.Range("S2").Value = 481
$ws.Range("T2").Value = 829
$ws.Range("U2").Value = 59
$ws.Range("V2").Value = 7003
$ws.Range("W2").Value = 3
$ws.Range("X2").Value = 6955
$ws.Range("Z2").Value = 111
$ws.Range("AA2").Value = 37
